$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new progress-report rows. Values are entered in the same
# order the original author must have typed them (this determines the
# order new entries are appended to the shared-strings table), which is
# why the column order per row looks a little unusual (name, end date,
# start date, task...).

# --- Person 1: Nguyen Thanh Huy ---
$ws.Range("A3").Value = "Nguyễn Thanh Huy"
$ws.Range("C3").Value = "24/12/2023"
$ws.Range("B3").Value = "25/11/2023"
$ws.Range("D3").Value = "Làm page đặt vé"
$ws.Range("D4").Value = "Chức năng thêm, xóa, sửa, xem user của trang admin, pagination "

# --- Person 2: Phan Vu Cong Thanh ---
$ws.Range("A6").Value = "Phan Vũ Công Thành"
$ws.Range("D6").Value = "Làm chức năng đặt vé, khi chọn ghế xong"
$ws.Range("D7").Value = "Chức năng thêm, xóa, sửa, xem danh sách phim "
$ws.Range("D8").Value = "Chức năng tạo lịch chiếu cho phim "

# Extra note added under person 1, entered last
$ws.Range("D5").Value = "Làm chức năng tìm kiếm user và phim"

# Reuse the same dates for person 2 (same shared-string values as row 3)
$ws.Range("B6").Value = "25/11/2023"
$ws.Range("C6").Value = "24/12/2023"

# Column D needed to widen considerably to fit the long task descriptions.
$ws.Columns.Item(4).ColumnWidth = 53

# Cursor ended up parked on E13 when the file was last saved.
[void]$ws.Range("E13").Select()
